# Split two paragraphs' single runs into multiple runs (no visible text-
# formatting change other than the inserted words), matching:
#   "I know that is has been..."   -> "I know that this has been..."
#     (split into: "I know that " | "th" | "is has been...")
#   "...go to school if I do..."   -> "...go to school even if I do..."
#     (split into: "...go to school " | "even " | "if I do...")
#
# Word's plain Range.InsertAfter/InsertBefore silently re-merges the
# touched run back into a single <w:r> on save (no distinguishing rsid),
# which would not reproduce the run split the diff shows. Recording the
# insert as a tracked change keeps the inserted text in its own run, and
# then accepting that specific revision (Revisions(i).Accept, not
# Document.AcceptAllRevisions which also strips sibling rsid attributes
# document-wide) leaves behind exactly the separate <w:r> elements the
# diff expects, with no leftover revision markup.

$d = $word.ActiveDocument

$d.TrackRevisions = $true

# --- Edit 1: "I know that is has been..." -> insert "th" before "is has" ---
$rng1 = $d.Content
$rng1.Find.Execute("I know that ", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$rng1.Collapse(0)
$rng1.InsertAfter("th")

# --- Edit 2: "...go to school if I do..." -> insert "even " before "if" ---
$rng2 = $d.Content
$rng2.Find.Execute("go to school ", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$rng2.Collapse(0)
$rng2.InsertAfter("even ")

$d.TrackRevisions = $false

# Accept just the revisions we created (leaves unrelated runs/rsids alone).
while ($d.Revisions.Count -gt 0) {
    $d.Revisions(1).Accept()
}
